$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the formulas/values from F11:F26 while keeping the cell formatting
$ws.Range("F11:F26").ClearContents()

# Update the active selection to match the author's last-selected cell
$ws.Range("F16").Select()
